$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Fill in the "Severity" column (D) for rows that were missing it
$ws.Range("D4").Value = "High"
$ws.Range("D5").Value = "Medium"
$ws.Range("D6").Value = "High"
$ws.Range("D7").Value = "High"
$ws.Range("D10").Value = "High"
$ws.Range("D12").Value = "Medium"
$ws.Range("D13").Value = "Medium"

# Fill in the "Assigned to" column (C) for rows that were missing it
$ws.Range("C14").Value = "Vinita"
$ws.Range("C15").Value = "Vinita"
$ws.Range("C16").Value = "Vinita"
$ws.Range("C17").Value = "Vinita"

# Update the active selection to reflect the last edited cell
$ws.Range("B15").Select()
